# Automatische test-sync: 2025-08-04 21:04:50
# Appends the new test-mail row (#18) to the "Logs" sheet, extends the
# conditional-formatting ranges to cover it, and bumps the corresponding
# category tally on the "Dashboard" sheet.

$wb = $excel.ActiveWorkbook

$logs = $wb.Worksheets.Item("Logs")
$dashboard = $wb.Worksheets.Item("Dashboard")

$newRow = 30

$logs.Cells.Item($newRow, 1).Value = "Bestel je 200 stuks M8-bouten RVS voor Van Dijk?"
$logs.Cells.Item($newRow, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item($newRow, 3).Value = "Testmail #18: Bestel je 200 stuks M8-bouten RVS voor Van Dijk?"
$logs.Cells.Item($newRow, 4).Value = "Inkoop / Bestellingen"
$logs.Cells.Item($newRow, 5).Value = "Bedankt, we hebben dit doorgestuurd naar inkoop@bedrijf.nl."
$logs.Cells.Item($newRow, 6).Value = "2025-08-04 21:04:43"
$logs.Cells.Item($newRow, 7).Value = "Ja"
$logs.Cells.Item($newRow, 8).Value = "Ja"
$logs.Cells.Item($newRow, 9).Value = "Nee"
$logs.Cells.Item($newRow, 10).Value = "Nee"

# Extend the existing conditional-formatting rules (columns D, G, H, I, J)
# so they keep covering the data down through the freshly added row 30,
# without disturbing their existing rules/dxf formats.
$logs.Range("D2:D29").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D30"))
$logs.Range("G2:G29").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G30"))
$logs.Range("H2:H29").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("H2:H30"))
$logs.Range("I2:I29").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("I2:I30"))
$logs.Range("J2:J29").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("J2:J30"))

# Bump the "Inkoop / Bestellingen" tally on the Dashboard sheet (6 -> 7).
$dashboard.Range("B3").Value = 7
